$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1): two new date columns (Jun_26, Jun_27) are inserted
# before the existing ones, pushing the old B:E header block to E:H.
#   B1 = Jun_27 (new)
#   C1 = Jun_26 (new)
#   D1 = Jun_26 (new, duplicate label)
#   E1 = Jun_17 (was B1)
#   F1 = Jun_15 (was C1)
#   G1 = Jun_13 (was D1)
#   H1 = Jun_10 (was E1)
# ---------------------------------------------------------------------------

$oldB1 = $ws.Range("B1").Value2
$oldC1 = $ws.Range("C1").Value2
$oldD1 = $ws.Range("D1").Value2
$oldE1 = $ws.Range("E1").Value2

$ws.Range("H1").Value = $oldE1
$ws.Range("G1").Value = $oldD1
$ws.Range("F1").Value = $oldC1
$ws.Range("E1").Value = $oldB1

$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# ---------------------------------------------------------------------------
# Data rows (2-27): the old "UN" rating column E holds either a plain "UN"
# or, for a few firms, a detailed rating-change note. That column slides to
# H, and the three newly-opened columns (E, F, G) are filled with "UN" to
# match the existing B:D columns.
# ---------------------------------------------------------------------------

for ($r = 2; $r -le 27; $r++) {
    $oldE = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 8).Value = $oldE
    $ws.Cells.Item($r, 5).Value = "UN"
    $ws.Cells.Item($r, 6).Value = "UN"
    $ws.Cells.Item($r, 7).Value = "UN"
}

# ---------------------------------------------------------------------------
# New rows for the newly added analyst/benchmark group.
# ---------------------------------------------------------------------------

$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"

# ---------------------------------------------------------------------------
# Column widths: mirror the existing custom width (8 chars) across the newly
# used columns F, G, H so the sheet keeps a uniform look.
# ---------------------------------------------------------------------------

$ws.Columns.Item(6).ColumnWidth = 7.14
$ws.Columns.Item(7).ColumnWidth = 7.14
$ws.Columns.Item(8).ColumnWidth = 7.14
